$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.226.11"
$ws.Range("E2").Value = "  -0.20%  "
$ws.Range("D3").Value = "2.638.54"
$ws.Range("E3").Value = "  -0.83%  "
$ws.Range("E4").Value = "  +0.02%  "
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.55"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  -2.37%  "
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.90"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  -2.81%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -1.26%  "
$ws.Range("D9").Value = "2.638.24"
$ws.Range("E9").Value = "  -0.82%  "
$ws.Range("E10").Value = "  -1.77%  "
$ws.Range("E12").Value = "  -0.97%  "
$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.94"
$ws.Range("D13").Style = $origStyle
$ws.Range("D14").Value = "3.132.05"
$ws.Range("E14").Value = "  -0.53%  "
$ws.Range("E15").Value = "  -2.18%  "
$ws.Range("D16").Value = "72.094.45"
$ws.Range("E16").Value = "  -0.27%  "
$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.84"
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = "  -2.80%  "
$ws.Range("D18").Value = "2.630.41"
$ws.Range("E18").Value = "  -1.34%  "
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.22"
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = "  +1.43%  "
$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.95"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = "  -0.74%  "
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "373.59"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  -1.68%  "
$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.15"
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = "  -1.37%  "
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.05"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  -1.35%  "
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("E25").Value = "  -2.38%  "
$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.26"
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = "  -3.23%  "
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.62"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = "  -4.12%  "
$ws.Range("D28").Value = "2.778.39"
$ws.Range("E28").Value = "  -1.02%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("E30").Value = "  -0.17%  "
$ws.Range("E31").Value = "  -2.86%  "
$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "495.74"
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = "  -5.04%  "
$ws.Range("E33").Value = "  -3.00%  "
$ws.Range("E34").Value = "  -1.24%  "
$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = "  +0.03%  "
$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "160.98"
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = "  -1.88%  "
$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.27"
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = "  -1.17%  "
$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.114"
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = "  +2.17%  "
$ws.Range("E39").Value = "  -1.22%  "
$ws.Range("E40").Value = "  -3.10%  "
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("E42").Value = "  -6.61%  "
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.55"
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = "  -2.69%  "
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.89"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = "  -4.12%  "
$ws.Range("E45").Value = "  -2.68%  "
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.10"
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = "  -0.68%  "
$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "152.70"
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = "  -0.18%  "
$ws.Range("B48").Value = "Filecoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.65"
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = "  -2.25%  "
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.547"
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = "  -0.28%  "
$ws.Range("E50").Value = "  -2.12%  "
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0747"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  -1.62%  "
